$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new log rows (287-292) to the bottom of the existing log table,
# matching the data added by the authoring tool for this commit.
$newRows = @(
    @("2025-12-08 05:50:56", "Admin", "Login", "login_success", "Role: admin"),
    @("2025-12-08 05:50:57", "Admin", "dashboard", "access_granted", "Opened dashboard page"),
    @("2025-12-08 05:50:59", "Admin", "settings", "access_granted", "Opened settings page"),
    @("2025-12-08 05:51:09", "Admin", "settings", "access_granted", "Opened settings page"),
    @("2025-12-08 05:51:13", "Admin", "settings", "access_granted", "Opened settings page"),
    @("2025-12-08 05:51:16", "Admin", "settings", "access_granted", "Opened settings page")
)

$startRow = 287
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
